$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New requirement rows (R3 already existed as a label, fill its details; then R4..R6 plus a trailing R7 label) ---

# R3 - Crear nuevas cuentas
$ws.Range("D6").Value = "Crear  nuevas cuentas"
$ws.Range("E6").Value = "permite crear nuevas cuentas si el usuario no tiene una"
$ws.Range("F6").Value = "datos del usuario"
$ws.Range("G6").Value = "nueva cuenta greada"

# R4 - Registrar alarmas
$ws.Range("D7").Value = "Registrar alarmas"
$ws.Range("E7").Value = "Le permite al usuario registrar las alarmas que quiera durante el dia "
$ws.Range("F7").Value = "Nombre de la alarma, tiempo que se ejecuta, "
$ws.Range("G7").Value = "Alarma guaradada"

# R5 - Mostrar alarma
$ws.Range("D8").Value = "Mostrar alarma"
$ws.Range("E8").Value = "Cuando el en el reloj del programa este a la hora de una alarma en especifico se mostrara un anuncio mostrando los datos de la alarma"
$ws.Range("F8").Value = "Tiempo que se ejecute la alarma"
$ws.Range("G8").Value = "Mensaje mostrando los datos de la alarma"

# R6 - Recomendar ejercicios
$ws.Range("D9").Value = "Recomendar ejercicios"
$ws.Range("E9").Value = "Se debe mostrar al usuario ejercios a diario con los datos del ejercicio"
$ws.Range("G9").Value = "Mostrar ejercicio"

# New trailing requirement label
$ws.Range("C10").Value = "R7"

# --- Row heights adjusted by Excel to fit the new/longer wrapped text ---
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 60
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 60
$ws.Rows.Item(9).RowHeight = 30

# --- Last selected cell in the sheet ---
$ws.Range("D10").Select()
